$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in E1, matching style of existing header cells (bold)
$ws.Range("E1").Value = "PriceExceedBy"
$ws.Range("E1").Font.Bold = $true

# Add new value in E2
$ws.Range("E2").Value = 200

# Best-fit the new column width to match the other bestFit header/data columns
$ws.Columns.Item(5).ColumnWidth = 12.333333333333334

# Move active selection to E1 as in the edited workbook
$ws.Range("E1").Select()
